$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new commit rows
$ws.Range("C28").Value = "base Runner & general update"
$ws.Range("G28").Value = 2.2

$ws.Range("C29").Value = "Runner update(80%) & general data update"
$ws.Range("G29").Value = 2.5

# Extend the total formula to include the new rows
$ws.Range("G39").Formula = "=SUM(G4:G29)"

# Update the selected cell to match the new active selection
$ws.Range("C30").Select()
